$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (copy style from an existing header cell, then set values)
$ws.Range("H1:H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Alpha"
$ws.Range("J1").Value = "Adjusted R2"

# Row 2 - model results updated (Newey-West), new Alpha/Adjusted R2 columns
$ws.Range("C2").Value = 0.1382729932387747
$ws.Range("E2").Value = -1.482252984112642
$ws.Range("I2").Value = 0.004877181779817691
$ws.Range("J2").Value = 0.2453894587614225

# Row 3
$ws.Range("C3").Value = 0.5359614640703878
$ws.Range("E3").Value = -0.6189315332752872
$ws.Range("I3").Value = -0.0007364963197267401
$ws.Range("J3").Value = 0.01316497873353151

# Row 4
$ws.Range("C4").Value = 0.243460710302973
$ws.Range("E4").Value = 1.166380323637777
$ws.Range("I4").Value = 0.002560966961518267
$ws.Range("J4").Value = 0.3764579905909564

# Row 5
$ws.Range("C5").Value = 0.5506088808635434
$ws.Range("E5").Value = 0.5968479800442473
$ws.Range("I5").Value = -0.002121081048118593
$ws.Range("J5").Value = 0.257457228890327

# Row 6 (Significant flips from TRUE to FALSE under the new model)
$ws.Range("C6").Value = 0.1597896828439402
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = -1.405779258898873
$ws.Range("I6").Value = 0.0112430145049063
$ws.Range("J6").Value = 0.1540426079179317

# Row 7
$ws.Range("C7").Value = 0.1858576202422962
$ws.Range("E7").Value = 1.322933119162968
$ws.Range("I7").Value = -0.002138288787201321
$ws.Range("J7").Value = 0.1556326995624852
